# Roll the GSC export window forward by one day:
# drop the oldest date row (2025-10-18) from the "Chart" sheet, shifting
# all remaining rows up by one, and mark the now-final two rows'
# "Impressions" values as 0 (finalized) instead of blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest day's row; everything below shifts up automatically.
$ws.Rows.Item(2).Delete()

# The last two rows used to be the "not yet finalized" rows (blank
# Impressions). After the shift, both now have confirmed 0 impressions.
$ws.Range("D89").Value = 0
$ws.Range("D90").Value = 0

Write-Output "done"
